# "Generate Report for Handback" — the localization-status report is
# regenerated: the de-de handback is now in sync with en-US (no more
# version-mismatch error), so the Status / Error Detail / handback
# timestamps for de-de (and the zh-cn status) are refreshed, and the
# columns that hold the longer status text are auto-fit wider while the
# (now empty) Error Detail column shrinks back down.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Columns.Item(5).ColumnWidth = 29.17
$ovw.Columns.Item(6).ColumnWidth = 29.17

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-06 05:03:32"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(16).ColumnWidth = 12.83

# --- de-de sheet ---------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-09-06 05:03:39"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(16).ColumnWidth = 12.83
